$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2339.7273
$ws.Range("I4").Value = 1962.4286
$ws.Range("J4").Value = 3000
$ws.Range("K4").Value = 1962.4286
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = -1848.4286
$ws.Range("N4").Value = -3228
$ws.Range("H70").Value = 6902.0835
$ws.Range("J70").Value = 8647
$ws.Range("L70").Value = 25941
$ws.Range("N70").Value = -26481
$ws.Range("H73").Value = 6902.0835
$ws.Range("J73").Value = 8647
$ws.Range("L73").Value = 25941
$ws.Range("N73").Value = -27813
$ws.Range("H86").Value = 4865.385
$ws.Range("I86").Value = 1900
$ws.Range("J86").Value = 5404.5454
$ws.Range("K86").Value = 1900
$ws.Range("L86").Value = 5404.5454
$ws.Range("M86").Value = -777
$ws.Range("N86").Value = -7650.5454
$ws.Range("H89").Value = 4865.385
$ws.Range("I89").Value = 1900
$ws.Range("J89").Value = 5404.5454
$ws.Range("K89").Value = 9500
$ws.Range("L89").Value = 27022.727
$ws.Range("M89").Value = -3884
$ws.Range("N89").Value = -38254.727
$ws.Range("H98").Value = 3273.125
$ws.Range("I98").Value = 1092.5
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 1092.5
$ws.Range("L98").Value = 4000
$ws.Range("M98").Value = 405.5
$ws.Range("N98").Value = -6996
$ws.Range("H122").Value = 3273.125
$ws.Range("I122").Value = 1092.5
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 3277.5
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -827.5
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 1128.4286
$ws.Range("I132").Value = 1067
$ws.Range("K132").Value = 3201
$ws.Range("M132").Value = -671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4554.1
$ws.Range("I61").Value = 2028.1666
$ws.Range("K61").Value = 2028.1666
$ws.Range("M61").Value = -1816.1666
$ws.Range("H82").Value = 42590
$ws.Range("J82").Value = 42590
$ws.Range("L82").Value = 42590
$ws.Range("N82").Value = -43312
$ws.Range("H85").Value = 42590
$ws.Range("J85").Value = 42590
$ws.Range("L85").Value = 42590
$ws.Range("N85").Value = -45086
$ws.Range("H122").Value = 1169.6666
$ws.Range("I122").Value = 1169.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3508.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1058.9998
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 51808
$ws.Range("J123").Value = 51808
$ws.Range("L123").Value = 51808
$ws.Range("N123").Value = -61608
$ws.Range("H132").Value = 2999
$ws.Range("I132").Value = 2665.3333
$ws.Range("K132").Value = 7995.999899999999
$ws.Range("M132").Value = -5465.999899999999
$ws.Range("H135").Value = 181999.42
$ws.Range("J135").Value = 181999.42
$ws.Range("L135").Value = 181999.42
$ws.Range("N135").Value = -192139.42
$ws.Range("H136").Value = 4554.1
$ws.Range("I136").Value = 2028.1666
$ws.Range("K136").Value = 6084.4998
$ws.Range("M136").Value = -3534.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 561.7
$ws.Range("I22").Value = 441.8889
$ws.Range("K22").Value = 441.8889
$ws.Range("M22").Value = -268.8889
$ws.Range("H92").Value = 38734
$ws.Range("J92").Value = 38734
$ws.Range("L92").Value = 38734
$ws.Range("N92").Value = -43726
$ws.Range("H95").Value = 17029.715
$ws.Range("J95").Value = 17029.715
$ws.Range("L95").Value = 17029.715
$ws.Range("N95").Value = -22521.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 130.04347
$ws.Range("J7").Value = 204.1
$ws.Range("L7").Value = 204.1
$ws.Range("N7").Value = -430.1
$ws.Range("H31").Value = 2589.8
$ws.Range("I31").Value = 2352
$ws.Range("J31").Value = 2827.6
$ws.Range("K31").Value = 2352
$ws.Range("L31").Value = 2827.6
$ws.Range("M31").Value = -2057
$ws.Range("N31").Value = -3417.6
$ws.Range("H34").Value = 2589.8
$ws.Range("I34").Value = 2352
$ws.Range("J34").Value = 2827.6
$ws.Range("K34").Value = 2352
$ws.Range("L34").Value = 2827.6
$ws.Range("M34").Value = -2150
$ws.Range("N34").Value = -3231.6
$ws.Range("H122").Value = 2904.5
$ws.Range("J122").Value = 2904.5
$ws.Range("L122").Value = 8713.5
$ws.Range("N122").Value = -13613.5
$ws.Range("H124").Value = 99999
$ws.Range("J124").Value = 99999
$ws.Range("L124").Value = 99999
$ws.Range("N124").Value = -104909
$ws.Range("H132").Value = 1404.6129
$ws.Range("I132").Value = 1255.1072
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 3765.3216
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -1235.3216
$ws.Range("N132").Value = -13460

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11111225
$ws.Range("I4").Value = 11111225
$ws.Range("K4").Value = 33333675
$ws.Range("M4").Value = -33333563
$ws.Range("H11").Value = 26804838
$ws.Range("I11").Value = 26804838
$ws.Range("K11").Value = 80414514
$ws.Range("M11").Value = -80414374
$ws.Range("H17").Value = 716.6667
$ws.Range("I17").Value = 75.5
$ws.Range("J17").Value = 1999
$ws.Range("K17").Value = 226.5
$ws.Range("L17").Value = 5997
$ws.Range("M17").Value = -57.5
$ws.Range("N17").Value = -6335
$ws.Range("H80").Value = 1401
$ws.Range("I80").Value = 1802
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 5406
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -4470
$ws.Range("N80").Value = -4872
$ws.Range("H83").Value = 1401
$ws.Range("I83").Value = 1802
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 16218
$ws.Range("L83").Value = 9000
$ws.Range("M83").Value = -11538
$ws.Range("N83").Value = -18360
$ws.Range("H113").Value = 1016
$ws.Range("I113").Value = 1026.5
$ws.Range("J113").Value = 1005.5
$ws.Range("K113").Value = 3079.5
$ws.Range("L113").Value = 3016.5
$ws.Range("M113").Value = -909.5
$ws.Range("N113").Value = -7356.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8368.9
$ws.Range("I70").Value = 8484.5
$ws.Range("K70").Value = 8484.5
$ws.Range("M70").Value = -8214.5
$ws.Range("H73").Value = 8368.9
$ws.Range("I73").Value = 8484.5
$ws.Range("K73").Value = 8484.5
$ws.Range("M73").Value = -7548.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1200
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 1200
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 1200
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -3696
$ws.Range("H97").Value = 18388.5
$ws.Range("J97").Value = 18388.5
$ws.Range("L97").Value = 18388.5
$ws.Range("N97").Value = -20370.5
$ws.Range("H136").Value = 2651.125
$ws.Range("I136").Value = 2223.7856
$ws.Range("J136").Value = 3249.4
$ws.Range("K136").Value = 6671.3568
$ws.Range("L136").Value = 9748.200000000001
$ws.Range("M136").Value = -4121.3568
$ws.Range("N136").Value = -14848.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14992.75
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 14992.75
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 14992.75
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -16240.75
$ws.Range("H65").Value = 14992.75
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 14992.75
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 74963.75
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -81203.75
$ws.Range("H126").Value = 5631.3335
$ws.Range("I126").Value = 3997.4
$ws.Range("J126").Value = 8899.200000000001
$ws.Range("K126").Value = 11992.2
$ws.Range("L126").Value = 26697.6
$ws.Range("M126").Value = -9522.200000000001
$ws.Range("N126").Value = -31637.6
$ws.Range("H136").Value = 2140.077
$ws.Range("I136").Value = 2276.7273
$ws.Range("J136").Value = 1388.5
$ws.Range("K136").Value = 6830.1819
$ws.Range("L136").Value = 4165.5
$ws.Range("M136").Value = -4280.1819
$ws.Range("N136").Value = -9265.5
